$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (even_MAG-GUT57690.fa), shifting
# it (and everything below it) down by one row.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new record.
$ws.Cells.Item(2, 1).Value = "even_MAG-GUT57658.fa"
$ws.Cells.Item(2, 2).Value = -16.04862608488929
$ws.Cells.Item(2, 3).Value = "s__CAG-791 sp000431495"
$ws.Cells.Item(2, 4).Value = "s__CAG-791 sp000431495(reject)"

# Match the formatting used on column A of the data rows (bold/border/
# centered alignment) without dragging that formatting onto columns B-D.
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Append a new record as row 4.
$ws.Cells.Item(4, 1).Value = "even_MAG-GUT58392.fa"
$ws.Cells.Item(4, 2).Value = -17.19533930795949
$ws.Cells.Item(4, 3).Value = "s__CAG-791 sp000431495"
$ws.Cells.Item(4, 4).Value = "s__CAG-791 sp000431495(reject)"

$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
